# Insert a new data row at row 641 (pushing the existing rows 641-693 down
# to 642-694) and populate the new row with a new price observation for
# "Feria Lagunitas de Puerto Montt" / Pomelo / Start Ruby / Primera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 641:693 down to 642:694, leaving row 641 empty for new data.
$ws.Rows(641).Insert()

# Populate the newly inserted row 641 with the new record.
$ws.Range("A641").Value = 4
$ws.Range("B641").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C641").Value = "Los Lagos"
$ws.Range("D641").Value = 45223
$ws.Range("E641").Value = 10
$ws.Range("F641").Value = "Fruta"
$ws.Range("G641").Value = 100102
$ws.Range("H641").Value = "Cítricos"
$ws.Range("I641").Value = 100102006
$ws.Range("J641").Value = "Pomelo"
$ws.Range("K641").Value = "Start Ruby"
$ws.Range("L641").Value = "Primera"
$ws.Range("M641").Value = 100
$ws.Range("N641").Value = 15000
$ws.Range("O641").Value = 15000
$ws.Range("P641").Value = 15000
$ws.Range("Q641").Value = "$/caja 14 kilos empedrada"
$ws.Range("R641").Value = "Región de O'Higgins"
$ws.Range("S641").Value = 1071
$ws.Range("T641").Value = 14
